$p = $ppt.ActivePresentation
$f = $p.Fonts
Write-Output ($f | Get-Member | Out-String)
Write-Output ("Count: " + $f.Count)
